# CoalesceOperator.xlsx edit
# - Adds a small "Account/Debit/Credit" table (_tExample) at B4:D8
# - Adds a second table (_tExample_) at G4:J6 that mirrors the first table
#   and adds a computed "Math" column (this is normally produced by a Power
#   Query "_tExample" query + query table/connection; the headless engine
#   here has no live Power Query / OLE DB engine, so we recreate the
#   resulting static shape: values, table objects, names, styles)
# - Adds the video-reference header row and a hidden ExternalData_1 name
# - Updates the selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / video reference row -----------------------------------------
$ws.Range("A1").Value = "Video Reference"
$ws.Range("C1").Value = "https://www.youtube.com/watch?v=-FQrxNMa_7I"

# --- First table: _tExample (B4:D8) ----------------------------------------
$ws.Range("B4").Value = "Account"
$ws.Range("C4").Value = "Debit"
$ws.Range("D4").Value = "Credit"

$ws.Range("B5").Value = "A"
$ws.Range("D5").Value = 2

$ws.Range("B6").Value = "B"
$ws.Range("C6").Value = 3

$ws.Range("B7").Value = "C"
$ws.Range("C7").Value = 4

$ws.Range("B8").Value = "D"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 1

# --- Second table: _tExample_ (G4:J6), the query-table result --------------
$ws.Range("G4").Value = "Account"
$ws.Range("H4").Value = "Debit"
$ws.Range("I4").Value = "Credit"
$ws.Range("J4").Value = "Math"

$ws.Range("G5").Value = "A"
$ws.Range("G5").NumberFormat = "General"
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = -2

$ws.Range("G6").Value = "B"
$ws.Range("G6").NumberFormat = "General"
$ws.Range("H6").Value = 3
$ws.Range("J6").Value = 3

# --- Create the two table (ListObject) structures ---------------------------
$null = $ws.ListObjects.Add(1, $ws.Range("B4:D8"), 0, 1)
$null = $ws.ListObjects.Add(1, $ws.Range("G4:J6"), 0, 1)

# Renaming a ListObject can reshuffle the (name-sorted) ListObjects
# collection, so re-resolve each table by its range address right before
# mutating it instead of caching object handles across calls.
function Get-ListObjectByAddress($sheet, $addr) {
    for ($i = 1; $i -le $sheet.ListObjects.Count; $i++) {
        $item = $sheet.ListObjects.Item($i)
        if ($item.Range.Address() -eq $addr) {
            return $item
        }
    }
    return $null
}

(Get-ListObjectByAddress $ws '$B$4:$D$8').Name = "_tExample"
(Get-ListObjectByAddress $ws '$B$4:$D$8').TableStyle = "Biegert Table Standard"

(Get-ListObjectByAddress $ws '$G$4:$J$6').Name = "_tExample_"
(Get-ListObjectByAddress $ws '$G$4:$J$6').TableStyle = "TableStyleMedium7"

# --- Hidden defined name pointing at the query table's result range --------
$extData = $ws.Names.Add("ExternalData_1", "=Report!`$G`$4:`$J`$6")
$extData.Visible = $false

# --- Column widths (best effort; engine quantizes these internally) --------
$ws.Columns.Item(2).ColumnWidth = 9.89
$ws.Columns.Item(7).ColumnWidth = 10.22
$ws.Columns.Item(8).ColumnWidth = 8.22
$ws.Columns.Item(9).ColumnWidth = 9.22
$ws.Columns.Item(10).ColumnWidth = 7.22

# --- Selection shown when the workbook is reopened --------------------------
$ws.Range("F14").Select() | Out-Null
